# Auto-generated edit script applying market-price data updates
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 6251730.5
$ws.Cells.Item(19, 10).Value = 9092836
$ws.Cells.Item(19, 12).Value = 9092836
$ws.Cells.Item(19, 14).Value = -9093186
$ws.Cells.Item(32, 8).Value = 4749.5
$ws.Cells.Item(32, 10).Value = 4749.5
$ws.Cells.Item(32, 12).Value = 4749.5
$ws.Cells.Item(32, 14).Value = -5401.5
$ws.Cells.Item(40, 8).Value = 100002240
$ws.Cells.Item(40, 9).Value = 1750
$ws.Cells.Item(40, 11).Value = 1750
$ws.Cells.Item(40, 13).Value = -1575
$ws.Cells.Item(88, 8).Value = 4541.143
$ws.Cells.Item(88, 10).Value = 4997.25
$ws.Cells.Item(88, 12).Value = 4997.25
$ws.Cells.Item(88, 14).Value = -5809.25
$ws.Cells.Item(91, 8).Value = 4541.143
$ws.Cells.Item(91, 10).Value = 4997.25
$ws.Cells.Item(91, 12).Value = 4997.25
$ws.Cells.Item(91, 14).Value = -7805.25
$ws.Cells.Item(97, 8).Value = 686.6667
$ws.Cells.Item(97, 10).Value = 686.6667
$ws.Cells.Item(97, 12).Value = 2060.0001
$ws.Cells.Item(97, 14).Value = -3052.0001
$ws.Cells.Item(98, 8).Value = 1416.3793
$ws.Cells.Item(98, 9).Value = 1216.5
$ws.Cells.Item(98, 11).Value = 1216.5
$ws.Cells.Item(98, 13).Value = 281.5
$ws.Cells.Item(122, 8).Value = 1416.3793
$ws.Cells.Item(122, 9).Value = 1216.5
$ws.Cells.Item(122, 11).Value = 3649.5
$ws.Cells.Item(122, 13).Value = -1199.5
$ws.Cells.Item(138, 8).Value = 4300.6553
$ws.Cells.Item(138, 9).Value = 2177.2917
$ws.Cells.Item(138, 11).Value = 6531.875100000001
$ws.Cells.Item(138, 13).Value = -1391.875100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(30, 8).Value = 45499
$ws.Cells.Item(30, 10).Value = 89989
$ws.Cells.Item(30, 12).Value = 89989
$ws.Cells.Item(30, 14).Value = -90289
$ws.Cells.Item(32, 8).Value = 2841.3333
$ws.Cells.Item(32, 9).Value = 2765.4666
$ws.Cells.Item(32, 11).Value = 2765.4666
$ws.Cells.Item(32, 13).Value = -2478.4666
$ws.Cells.Item(45, 8).Value = 3243.5557
$ws.Cells.Item(45, 9).Value = 1243.3334
$ws.Cells.Item(45, 11).Value = 1243.3334
$ws.Cells.Item(45, 13).Value = -866.3334
$ws.Cells.Item(61, 8).Value = 35002816
$ws.Cells.Item(61, 9).Value = 40003380
$ws.Cells.Item(61, 11).Value = 40003380
$ws.Cells.Item(61, 13).Value = -40003168
$ws.Cells.Item(102, 8).Value = 2842.2666
$ws.Cells.Item(102, 9).Value = 2358
$ws.Cells.Item(102, 10).Value = 5990
$ws.Cells.Item(102, 11).Value = 2358
$ws.Cells.Item(102, 12).Value = 5990
$ws.Cells.Item(102, 13).Value = -736
$ws.Cells.Item(102, 14).Value = -9234
$ws.Cells.Item(132, 8).Value = 4767429
$ws.Cells.Item(132, 9).Value = 5852.8423
$ws.Cells.Item(132, 11).Value = 17558.5269
$ws.Cells.Item(132, 13).Value = -15028.5269
$ws.Cells.Item(136, 8).Value = 35002816
$ws.Cells.Item(136, 9).Value = 40003380
$ws.Cells.Item(136, 11).Value = 120010140
$ws.Cells.Item(136, 13).Value = -120007590

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 8632.583000000001
$ws.Cells.Item(20, 10).Value = 2408.8572
$ws.Cells.Item(20, 12).Value = 2408.8572
$ws.Cells.Item(20, 14).Value = -2902.8572
$ws.Cells.Item(86, 8).Value = 2813.7407
$ws.Cells.Item(86, 9).Value = 2040.3334
$ws.Cells.Item(86, 10).Value = 3780.5
$ws.Cells.Item(86, 11).Value = 2040.3334
$ws.Cells.Item(86, 12).Value = 3780.5
$ws.Cells.Item(86, 13).Value = -917.3334
$ws.Cells.Item(86, 14).Value = -6026.5
$ws.Cells.Item(89, 8).Value = 2813.7407
$ws.Cells.Item(89, 9).Value = 2040.3334
$ws.Cells.Item(89, 10).Value = 3780.5
$ws.Cells.Item(89, 11).Value = 10201.667
$ws.Cells.Item(89, 12).Value = 18902.5
$ws.Cells.Item(89, 13).Value = -4585.666999999999
$ws.Cells.Item(89, 14).Value = -30134.5
$ws.Cells.Item(94, 8).Value = 1753.6
$ws.Cells.Item(94, 9).Value = 1423.7059
$ws.Cells.Item(94, 11).Value = 1423.7059
$ws.Cells.Item(94, 13).Value = -972.7058999999999
$ws.Cells.Item(105, 8).Value = 1147140.1
$ws.Cells.Item(105, 9).Value = 1526633.6
$ws.Cells.Item(105, 10).Value = 8659.6
$ws.Cells.Item(105, 11).Value = 1526633.6
$ws.Cells.Item(105, 12).Value = 8659.6
$ws.Cells.Item(105, 13).Value = -1524886.6
$ws.Cells.Item(105, 14).Value = -12153.6
$ws.Cells.Item(107, 8).Value = 4002.3333
$ws.Cells.Item(107, 10).Value = 6074.75
$ws.Cells.Item(107, 12).Value = 6074.75
$ws.Cells.Item(107, 14).Value = -9914.75
$ws.Cells.Item(134, 8).Value = 5558312.5
$ws.Cells.Item(134, 9).Value = 2175.1428
$ws.Cells.Item(134, 10).Value = 9094036
$ws.Cells.Item(134, 11).Value = 6525.428400000001
$ws.Cells.Item(134, 12).Value = 27282108
$ws.Cells.Item(134, 13).Value = -3990.428400000001
$ws.Cells.Item(134, 14).Value = -27287178

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 48108840
$ws.Cells.Item(31, 9).Value = 58827824
$ws.Cells.Item(31, 11).Value = 58827824
$ws.Cells.Item(31, 13).Value = -58827529
$ws.Cells.Item(34, 8).Value = 48108840
$ws.Cells.Item(34, 9).Value = 58827824
$ws.Cells.Item(34, 11).Value = 58827824
$ws.Cells.Item(34, 13).Value = -58827622
$ws.Cells.Item(52, 8).Value = 67999.5
$ws.Cells.Item(52, 10).Value = 99999
$ws.Cells.Item(52, 12).Value = 99999
$ws.Cells.Item(52, 14).Value = -100587
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).Value = ""
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).Value = ""
$ws.Cells.Item(94, 8).Value = 1121.5
$ws.Cells.Item(94, 10).Value = 1072.5555
$ws.Cells.Item(94, 12).Value = 1072.5555
$ws.Cells.Item(94, 14).Value = -1974.5555
$ws.Cells.Item(122, 8).Value = 3801.5
$ws.Cells.Item(122, 9).Value = 3728.7144
$ws.Cells.Item(122, 11).Value = 11186.1432
$ws.Cells.Item(122, 13).Value = -8736.143199999999
$ws.Cells.Item(132, 8).Value = 2587
$ws.Cells.Item(132, 9).Value = 2297.52
$ws.Cells.Item(132, 10).Value = 3244.9092
$ws.Cells.Item(132, 11).Value = 6892.559999999999
$ws.Cells.Item(132, 12).Value = 9734.7276
$ws.Cells.Item(132, 13).Value = -4362.559999999999
$ws.Cells.Item(132, 14).Value = -14794.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 7439
$ws.Cells.Item(11, 9).Value = 2450.25
$ws.Cells.Item(11, 10).Value = 17416.5
$ws.Cells.Item(11, 11).Value = 7350.75
$ws.Cells.Item(11, 12).Value = 52249.5
$ws.Cells.Item(11, 13).Value = -7210.75
$ws.Cells.Item(11, 14).Value = -52529.5
$ws.Cells.Item(14, 8).Value = 18042
$ws.Cells.Item(14, 9).Value = 18042
$ws.Cells.Item(14, 11).Value = 54126
$ws.Cells.Item(14, 13).Value = -53953
$ws.Cells.Item(23, 8).Value = 1443.4546
$ws.Cells.Item(23, 9).Value = 99
$ws.Cells.Item(23, 10).Value = 1577.9
$ws.Cells.Item(23, 11).Value = 297
$ws.Cells.Item(23, 12).Value = 4733.700000000001
$ws.Cells.Item(23, 13).Value = -62
$ws.Cells.Item(23, 14).Value = -5203.700000000001
$ws.Cells.Item(122, 8).Value = 30414.545
$ws.Cells.Item(122, 10).Value = 410.4
$ws.Cells.Item(122, 12).Value = 3693.6
$ws.Cells.Item(122, 14).Value = -8593.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2501.2666
$ws.Cells.Item(102, 9).Value = 2515.7144
$ws.Cells.Item(102, 11).Value = 2515.7144
$ws.Cells.Item(102, 13).Value = -893.7143999999998
$ws.Cells.Item(132, 8).Value = 7957676.5
$ws.Cells.Item(132, 9).Value = 3397.3333
$ws.Cells.Item(132, 10).Value = 21214808
$ws.Cells.Item(132, 11).Value = 10191.9999
$ws.Cells.Item(132, 12).Value = 63644424
$ws.Cells.Item(132, 13).Value = -7661.999899999999
$ws.Cells.Item(132, 14).Value = -63649484

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5919.4546
$ws.Cells.Item(7, 9).Value = 5511.5
$ws.Cells.Item(7, 11).Value = 5511.5
$ws.Cells.Item(7, 13).Value = -5399.5
$ws.Cells.Item(100, 8).Value = 20858250
$ws.Cells.Item(100, 9).Value = 3972.5
$ws.Cells.Item(100, 11).Value = 3972.5
$ws.Cells.Item(100, 13).Value = -3431.5
$ws.Cells.Item(126, 8).Value = 5919.4546
$ws.Cells.Item(126, 9).Value = 5511.5
$ws.Cells.Item(126, 11).Value = 16534.5
$ws.Cells.Item(126, 13).Value = -14064.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1002469
$ws.Cells.Item(132, 9).Value = 2098.5715
$ws.Cells.Item(132, 11).Value = 6295.7145
$ws.Cells.Item(132, 13).Value = -3765.7145
$ws.Cells.Item(136, 8).Value = 456429.1
$ws.Cells.Item(136, 9).Value = 1921.75
$ws.Cells.Item(136, 11).Value = 5765.25
$ws.Cells.Item(136, 13).Value = -3215.25
